$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "317.74"
$ws.Range("E2").Value = "-3.47%"
$ws.Range("D3").Value = "41.59"
$ws.Range("E3").Value = "-5.24%"
$ws.Range("D4").Value = "5.186"
$ws.Range("E4").Value = "-1.79%"
$ws.Range("D5").Value = "0.08107"
$ws.Range("E5").Value = "-3.46%"
$ws.Range("D6").Value = "4.367"
$ws.Range("E6").Value = "-1.61%"
$ws.Range("E7").Value = "-10.38%"
$ws.Range("D8").Value = "0.9300"
$ws.Range("E8").Value = "-4.67%"
$ws.Range("D9").Value = "0.1121"
$ws.Range("E9").Value = "-1.77%"
$ws.Range("D10").Value = "0.1852"
$ws.Range("E10").Value = "-2.53%"
$ws.Range("D11").Value = "0.09217"
$ws.Range("E11").Value = "-4.58%"
$ws.Range("D12").Value = "0.04587"
$ws.Range("E12").Value = "-0.60%"
$ws.Range("D13").Value = "7.389"
$ws.Range("E13").Value = "-18.18%"
$ws.Range("D14").Value = "0.1056"
$ws.Range("E14").Value = "-0.64%"
$ws.Range("D15").Value = "0.001276"
$ws.Range("E15").Value = "-2.35%"
$ws.Range("D16").Value = "0.005840"
$ws.Range("E16").Value = "-3.89%"
$ws.Range("D17").Value = "3.345"
$ws.Range("E17").Value = "-1.68%"
$ws.Range("E18").Value = "3.65%"
$ws.Range("D19").Value = "0.3369"
$ws.Range("E19").Value = "1.19%"
$ws.Range("D20").Value = "0.1385"
$ws.Range("E20").Value = "0.96%"
$ws.Range("D22").Value = "0.04201"
$ws.Range("E22").Value = "0.87%"
$ws.Range("E23").Value = "-3.95%"
$ws.Range("D24").Value = "0.004260"
$ws.Range("E24").Value = "-3.37%"
$ws.Range("D25").Value = "0.0001223"
$ws.Range("E25").Value = "-5.98%"
$ws.Range("D26").Value = "0.0002991"
$ws.Range("E26").Value = "0.14%"
$ws.Range("D38").Value = "0.02578"
$ws.Range("E38").Value = "-3.20%"
$ws.Range("D39").Value = "0.05478"
$ws.Range("E39").Value = "-2.88%"
$ws.Range("D40").Value = "0.008066"
$ws.Range("E40").Value = "2.56%"
$ws.Range("D41").Value = "0.1391"
$ws.Range("E41").Value = "-1.67%"
$ws.Range("D42").Value = "0.006554"
$ws.Range("E42").Value = "-10.94%"
$ws.Range("D43").Value = "0.002093"
$ws.Range("E43").Value = "0.76%"
$ws.Range("D44").Value = "0.008240"
$ws.Range("E44").Value = "4.31%"
$ws.Range("D45").Value = "0.3456"
$ws.Range("E45").Value = "-1.72%"
$ws.Range("D46").Value = "0.00006762"
$ws.Range("E46").Value = "-1.97%"
$ws.Range("E47").Value = "0.19%"
$ws.Range("D48").Value = "0.003383"
$ws.Range("E48").Value = "-3.49%"
$ws.Range("D49").Value = "0.004117"
$ws.Range("E49").Value = "16.54%"
$ws.Range("D50").Value = "0.00002106"
$ws.Range("E50").Value = "0.19%"
$ws.Range("D51").Value = "0.0002005"
$ws.Range("E51").Value = "0.19%"
